# Auto-generated Excel COM-interop script applying the cryptos.xlsx price/volume update
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "23.686.99"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.35%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.654.27"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.22%  "

$ws.Range("E4").Value = "  -0.16%  "

$ws.Range("E5").Value = "  -0.05%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "303.07"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.05%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3841"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.37%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3604"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.92%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "50.99"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.11%  "

$ws.Range("E10").Value = "  +0.28%  "

$ws.Range("E11").Value = "  +0.83%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.001"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.23%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "22.41"
$ws.Range("D13").Style = "Normal"

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.458"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.60%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.439"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.94%  "

$ws.Range("E16").Value = "  -0.86%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.653.44"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.72%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "97.57"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.68%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.07033"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.17%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.769"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +3.09%  "

$ws.Range("E21").Value = "  +1.23%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.002"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.13%  "

$ws.Range("E23").Value = "  +1.74%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "23.684.09"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.32%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.473"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -3.23%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.029"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.60%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "21.24"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.87%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "154.02"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.45%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.246"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.39%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "134.07"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.44%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.836.78"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.51%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.085"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +9.04%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.253"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +4.87%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "12.08"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +5.89%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.059"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.17%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02807"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.90%  "

$ws.Range("E37").Value = "  +0.36%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.08796"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.40%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "6.077"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.29%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.06982"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.54%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "13.03"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +6.82%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.6986"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.02%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.337"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.63%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "15.99"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +3.39%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.6504"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.91%  "

$ws.Range("E46").Value = "  -0.09%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.301"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.10%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.963"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.12%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.07873"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.78%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "128.19"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.30%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.179"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.62%  "
